# VyTrackLoginPage DDT: add a "VyTrackQA2User" sheet with credentials/results
# read & written by VyTrackLoginDDTTest.java, and normalize the Employees
# sheet's Salary column (drop the stray float typing on E2:E5).

$wb = $excel.ActiveWorkbook
$wsEmployees = $wb.Worksheets.Item("Employees")

# --- New worksheet, placed right after "Employees" ---
$ws2 = $wb.Worksheets.Add($null, $wsEmployees)
$ws2.Name = "VyTrackQA2User"

# Header row
$headers = @("username", "password", "firstname", "lastname", "Result")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws2.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Data rows: username, password, firstname, lastname, Result
$ws2.Cells.Item(2, 1).Value = "user1"
$ws2.Cells.Item(2, 2).Value = "UserUser123"
$ws2.Cells.Item(2, 3).Value = "John"
$ws2.Cells.Item(2, 4).Value = "Doe"
$ws2.Cells.Item(2, 5).Value = "PASSED"

$ws2.Cells.Item(3, 1).Value = "user2"
$ws2.Cells.Item(3, 2).Value = "UserUser123"
$ws2.Cells.Item(3, 3).Value = "Bella"
$ws2.Cells.Item(3, 4).Value = "Stamm"
$ws2.Cells.Item(3, 5).Value = "PASSED"

$ws2.Cells.Item(4, 1).Value = "storemanager51"
$ws2.Cells.Item(4, 2).Value = "UserUser123"
$ws2.Cells.Item(4, 3).Value = "Edd"
$ws2.Cells.Item(4, 4).Value = "Turner"
$ws2.Cells.Item(4, 5).Value = "PASSED"

$ws2.Cells.Item(5, 2).Value = "UserUser123"
$ws2.Cells.Item(5, 3).Value = "Roma"
$ws2.Cells.Item(5, 4).Value = "Medhurst"
$ws2.Cells.Item(5, 1).Value = "storemanager52"
$ws2.Cells.Item(5, 5).Value = "PASSED"

$ws2.Cells.Item(6, 1).Value = "storemanager101"
$ws2.Cells.Item(6, 2).Value = "UserUser123"
$ws2.Cells.Item(6, 3).Value = "John"
$ws2.Cells.Item(6, 4).Value = "Doe"
$ws2.Cells.Item(6, 5).Value = "PASSED"

$ws2.Cells.Item(7, 1).Value = "storemanager102"
$ws2.Cells.Item(7, 2).Value = "UserUser123"
$ws2.Cells.Item(7, 3).Value = "John"
$ws2.Cells.Item(7, 4).Value = "Doe"
$ws2.Cells.Item(7, 5).Value = "PASSED"

# Column widths (best-fit-like custom widths for the username/password columns)
$ws2.Columns.Item(1).ColumnWidth = 13.830729166666666
$ws2.Columns.Item(2).ColumnWidth = 10.998697916666666

# View: the new sheet becomes the active/selected tab, zoomed in, with C6:D7 selected
$ws2.Activate()
$ws2.Range("C6:D7").Select()
$excel.ActiveWindow.Zoom = 218

# Employees sheet: Salary column values should no longer be tagged as an
# explicit numeric type / float literal -- rewrite them as plain numbers.
$wsEmployees.Cells.Item(2, 5).Value = 200000
$wsEmployees.Cells.Item(3, 5).Value = 190000
$wsEmployees.Cells.Item(4, 5).Value = 175000
$wsEmployees.Cells.Item(5, 5).Value = 210000

# Keep Employees' own selection as it was (F1:F5); the newly active sheet is VyTrackQA2User.
$wsEmployees.Range("F1:F5").Select()
$ws2.Activate()
$ws2.Range("C6:D7").Select()
